$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-18 20:38:09"
$wsZhCn.Range("H4").Value = "2016-03-18 20:38:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-18 20:38:12"
$wsDeDe.Range("H4").Value = "2016-03-18 20:38:34"
